# Adds a new "3. Marketing and Public Relations" section to the document,
# right after the "2. Gaming and Competitive Teams" section and before the
# trailing blank paragraph(s) at the end of the body.

$d = $word.ActiveDocument

# Locate the last real paragraph of section 2 ("Compete in tournaments...")
# so we don't depend on hard-coded paragraph indices.
$anchor = $d.Content
$found = $anchor.Find.Execute("Compete in tournaments and represent the team")
if (-not $found) {
    throw "Anchor paragraph for section 2 not found"
}

# Expand the found text to the whole paragraph (includes the paragraph mark),
# so .End lands exactly at the start of the following paragraph.
$anchor.Expand(4) | Out-Null

# New section content expressed as WordprocessingML, wrapped in the
# single-part package format InsertXML expects.
$newSectionXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:lastRenderedPageBreak/>
              <w:t>3. Marketing and Public Relations</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:tab/>
              <w:t>&#8226;</w:t>
            </w:r>
            <w:r>
              <w:tab/>
              <w:t>Responsibilities:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>Promote the brand, manage social media platforms, create partnerships, and handle public relations.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:tab/>
              <w:t>&#8226;</w:t>
            </w:r>
            <w:r>
              <w:tab/>
              <w:t>Marketing Manager</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">Making a good product to attract the attention of people to became supporters of team such as </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>mafla</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>, jersey team, handicap and so on.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:tab/>
              <w:t>&#8226;</w:t>
            </w:r>
            <w:r>
              <w:tab/>
              <w:t>Content Creators</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>Scout the famous person to became brand ambassador to promote the team in public and social media such as TikTok, YouTube Facebook and others.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

# Insert right at the boundary between section 2's last paragraph and the
# blank paragraph that follows it, so the existing trailing blank paragraph
# is preserved as-is after the newly inserted content.
$insertionPoint = $d.Range($anchor.End, $anchor.End)
$insertionPoint.InsertXML($newSectionXml)
